$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing row (328) down through the new rows (329-343)
$ws.Range("A328:D328").Copy() | Out-Null
$ws.Range("A329:D343").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New daily data rows (dates 2021-07-26 through 2021-08-09)
$data = @(
    ,@(329, 44403, 0, 1, 15.22997258604935)
    ,@(330, 44404, 0, 1, 15.22997258604935)
    ,@(331, 44405, 0, 1, 15.22997258604935)
    ,@(332, 44406, 1, 1, 15.22997258604935)
    ,@(333, 44407, 0, 1, 15.22997258604935)
    ,@(334, 44408, 0, 1, 15.22997258604935)
    ,@(335, 44409, 2, 3, 45.68991775814803)
    ,@(336, 44410, 2, 5, 76.14986293024673)
    ,@(337, 44411, 1, 6, 91.37983551629607)
    ,@(338, 44412, 1, 7, 106.6098081023454)
    ,@(339, 44413, 1, 7, 106.6098081023454)
    ,@(340, 44414, 0, 7, 106.6098081023454)
    ,@(341, 44415, 2, 9, 137.0697532744441)
    ,@(342, 44416, 0, 7, 106.6098081023454)
    ,@(343, 44417, 1, 6, 91.37983551629607)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
}
